$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (entire rows), keeping rows 1-4
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2 updates
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 3.222529
$ws.Range("N2").Value = 9.667587000000001
$ws.Range("O2").Value = 0.1507607113527322
$ws.Range("P2").Value = 0.1507607113527322
$ws.Range("Q2").Value = 5.170463994746
$ws.Range("R2").Value = 46.53417595271401
$ws.Range("S2").Value = 0.1507607113527322
$ws.Range("T2").Value = 0.1507607113527322

# Row 3 updates
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 0.5693303333333334
$ws.Range("N3").Value = 1.707991
$ws.Range("O3").Value = 0.02663518188603469
$ws.Range("P3").Value = 0.02663518188603469
$ws.Range("Q3").Value = 0.9134757172446667
$ws.Range("R3").Value = 8.221281455202
$ws.Range("S3").Value = 0.02663518188603469
$ws.Range("T3").Value = 0.02663518188603469

# Row 4 updates
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 17.58326533333333
$ws.Range("N4").Value = 52.749796
$ws.Range("O4").Value = 0.8226041067612331
$ws.Range("P4").Value = 0.8226041067612331
$ws.Range("Q4").Value = 28.21189206243466
$ws.Range("R4").Value = 253.907028561912
$ws.Range("S4").Value = 0.8226041067612331
$ws.Range("T4").Value = 0.8226041067612331
